# Flights.xlsx — "Aircraft functionality added, and the import of the
# importFlights added."
#
# 1. Tidy-up of rows 4 and 5's flight numbers: they become plain integers
#    (3 / 4) instead of the stray 3.0 / 4.0 the earlier import left behind.
# 2. Five new flights (rows 6-10) are imported — this is the
#    "importFlights" feature the commit message refers to. Row 10 reuses
#    an existing "Taca"/"rf56"-style flight mixed with new destinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. normalize the stray 3.0 / 4.0 left on rows 4 and 5 -----------------
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# --- 2. import five new flights into rows 6-10 ------------------------------
# Column F holds dates that look like "M/D/YY" or "MM/DD/YYYY" — left alone,
# Excel's smart entry would silently convert those into date serials. The
# imported rows must stay as literal text, so each date cell is forced to
# the Text number format right before the value is written, then the
# cell style is put back to Normal afterwards (this is exactly what the
# "importFlights" text-import routine that produced rows 4/5 already did).

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 6
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Avianca"
$ws.Range("C6").Value = "rf56"
$ws.Range("D6").Value = "san salvador"
$ws.Range("E6").Value = "hawaii"
Set-TextValue $ws.Range("F6") "6/10/21"
$ws.Range("G6").Value = "2:00"
$ws.Range("H6").Value = "23:00"

# Row 7
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Taca"
$ws.Range("C7").Value = "hh7"
$ws.Range("D7").Value = "Paris"
$ws.Range("E7").Value = "new york"
Set-TextValue $ws.Range("F7") "6/12/21"
$ws.Range("G7").Value = "7:00"
$ws.Range("H7").Value = "19:00"

# Row 8
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "United"
$ws.Range("C8").Value = "SFDS3"
$ws.Range("D8").Value = "Toronto"
$ws.Range("E8").Value = "California"
Set-TextValue $ws.Range("F8") "05/02/2021"
$ws.Range("G8").Value = "10:00"
$ws.Range("H8").Value = "2:00"

# Row 9
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Delta"
$ws.Range("C9").Value = "SDFD4"
$ws.Range("D9").Value = "milan"
$ws.Range("E9").Value = "sevillaa"
Set-TextValue $ws.Range("F9") "12/05/2021"
$ws.Range("G9").Value = "15:00"
$ws.Range("H9").Value = "23:00"

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Taca"
$ws.Range("C10").Value = "RF52"
$ws.Range("D10").Value = "San Salvador"
$ws.Range("E10").Value = "Lima, Peru"
Set-TextValue $ws.Range("F10") "14/06/2021"
$ws.Range("G10").Value = "5:50"
$ws.Range("H10").Value = "10:00"

# --- selection mirrors where the author's cursor ended up after pasting ----
$ws.Range("A6:H6").Select() | Out-Null
